$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.443613805178813
$ws.Cells.Item(2, 3).Value = 0.5508010229219167
$ws.Cells.Item(2, 4).Value = 0.08985838777794442
$ws.Cells.Item(2, 6).Value = 2.597802277714038
$ws.Cells.Item(2, 7).Value = 0.002515046553118214
$ws.Cells.Item(2, 9).Value = 1.183514997888679
$ws.Cells.Item(2, 10).Value = 0.2112062533702854
$ws.Cells.Item(2, 13).Value = 0.6031615746933383
$ws.Cells.Item(2, 14).Value = 1.731397487175229

$ws.Cells.Item(3, 2).Value = 1.340194044334282
$ws.Cells.Item(3, 3).Value = 0.511404890779886
$ws.Cells.Item(3, 4).Value = 0.08968649554753583
$ws.Cells.Item(3, 6).Value = 2.571533506554303
$ws.Cells.Item(3, 7).Value = 0.002520348434799307
$ws.Cells.Item(3, 9).Value = 1.178859154891818
$ws.Cells.Item(3, 10).Value = 0.2111375892942249
$ws.Cells.Item(3, 13).Value = 0.5769380355127041
$ws.Cells.Item(3, 14).Value = 1.751163945221677

$ws.Cells.Item(4, 2).Value = 1.27745317462734
$ws.Cells.Item(4, 3).Value = 0.4875346308964481
$ws.Cells.Item(4, 4).Value = 0.08959376225058335
$ws.Cells.Item(4, 6).Value = 2.55699872763924
$ws.Cells.Item(4, 7).Value = 0.002523774712830032
$ws.Cells.Item(4, 9).Value = 1.176722879374132
$ws.Cells.Item(4, 10).Value = 0.2111891970457691
$ws.Cells.Item(4, 13).Value = 0.5611828721971435
$ws.Cells.Item(4, 14).Value = 1.763951319431428

$ws.Cells.Item(5, 2).Value = 1.25207568098898
$ws.Cells.Item(5, 3).Value = 0.4778868110935548
$ws.Cells.Item(5, 4).Value = 0.08955921180413284
$ws.Cells.Item(5, 6).Value = 2.551475053909257
$ws.Cells.Item(5, 7).Value = 0.002525214071893922
$ws.Cells.Item(5, 9).Value = 1.176033263448275
$ws.Cells.Item(5, 10).Value = 0.2112337804331545
$ws.Cells.Item(5, 13).Value = 0.5548492861010317
$ws.Cells.Item(5, 14).Value = 1.769325527606643

$ws.Cells.Item(6, 2).Value = 1.247873209404759
$ws.Cells.Item(6, 3).Value = 0.4762895784399461
$ws.Cells.Item(6, 4).Value = 0.08955367090288391
$ws.Cells.Item(6, 6).Value = 2.550581928032003
$ws.Cells.Item(6, 7).Value = 0.00252545568530518
$ws.Cells.Item(6, 9).Value = 1.175929660066735
$ws.Cells.Item(6, 10).Value = 0.2112426052649141
$ws.Cells.Item(6, 13).Value = 0.5538028316260579
$ws.Cells.Item(6, 14).Value = 1.770227759852752

$ws.Cells.Item(7, 2).Value = 1.277110157213997
$ws.Cells.Item(7, 3).Value = 0.4874041960808313
$ws.Cells.Item(7, 4).Value = 0.08959328315092563
$ws.Cells.Item(7, 6).Value = 2.556922618483014
$ws.Cells.Item(7, 7).Value = 0.002523793949900964
$ws.Cells.Item(7, 9).Value = 1.176712847352505
$ws.Cells.Item(7, 10).Value = 0.2111897029820824
$ws.Cells.Item(7, 13).Value = 0.5610971042775645
$ws.Cells.Item(7, 14).Value = 1.764023137401768

$ws.Cells.Item(8, 2).Value = 1.40779643464316
$ws.Cells.Item(8, 3).Value = 0.5371505400428873
$ws.Cells.Item(8, 4).Value = 0.08979647144503033
$ws.Cells.Item(8, 6).Value = 2.588412946772337
$ws.Cells.Item(8, 7).Value = 0.002516839260849224
$ws.Cells.Item(8, 9).Value = 1.181759278345851
$ws.Cells.Item(8, 10).Value = 0.2111631001949306
$ws.Cells.Item(8, 13).Value = 0.5940476634124678
$ws.Cells.Item(8, 14).Value = 1.738077688802072

$ws.Cells.Item(9, 2).Value = 1.670157353677496
$ws.Cells.Item(9, 3).Value = 0.637274322706503
$ws.Cells.Item(9, 4).Value = 0.09029582579662332
$ws.Cells.Item(9, 6).Value = 2.662892806538508
$ws.Cells.Item(9, 7).Value = 0.002504550385137179
$ws.Cells.Item(9, 9).Value = 1.197422566861363
$ws.Cells.Item(9, 10).Value = 0.2118563729647249
$ws.Cells.Item(9, 13).Value = 0.6614276517381086
$ws.Cells.Item(9, 14).Value = 1.692377655767622

$ws.Cells.Item(10, 2).Value = 1.866730667961178
$ws.Cells.Item(10, 3).Value = 0.7124655381932143
$ws.Cells.Item(10, 4).Value = 0.09072332424761242
$ws.Cells.Item(10, 6).Value = 2.725489098771618
$ws.Cells.Item(10, 7).Value = 0.002496334782802112
$ws.Cells.Item(10, 9).Value = 1.21249813365921
$ws.Cells.Item(10, 10).Value = 0.2128225400049431
$ws.Cells.Item(10, 13).Value = 0.7126468111133732
$ws.Cells.Item(10, 14).Value = 1.661980067121675

$ws.Cells.Item(11, 2).Value = 1.957010555473175
$ws.Cells.Item(11, 3).Value = 0.7470401772769719
$ws.Cells.Item(11, 4).Value = 0.09093079060955844
$ws.Cells.Item(11, 6).Value = 2.755702183580979
$ws.Cells.Item(11, 7).Value = 0.002492771805247574
$ws.Cells.Item(11, 9).Value = 1.22014264002361
$ws.Cells.Item(11, 10).Value = 0.2133618234464691
$ws.Cells.Item(11, 13).Value = 0.7363271480739115
$ws.Cells.Item(11, 14).Value = 1.648846144778972

$ws.Cells.Item(12, 2).Value = 1.991322133945062
$ws.Cells.Item(12, 3).Value = 0.7601868689741309
$ws.Cells.Item(12, 4).Value = 0.09101120567502718
$ws.Cells.Item(12, 6).Value = 2.767394954643407
$ws.Cells.Item(12, 7).Value = 0.002491447513188691
$ws.Cells.Item(12, 9).Value = 1.223151419603028
$ws.Cells.Item(12, 10).Value = 0.2135804234291854
$ws.Cells.Item(12, 13).Value = 0.7453494442229101
$ws.Cells.Item(12, 14).Value = 1.643972912215148

$ws.Cells.Item(13, 2).Value = 1.983926962857993
$ws.Cells.Item(13, 3).Value = 0.7573530771182959
$ws.Cells.Item(13, 4).Value = 0.09099380474263796
$ws.Cells.Item(13, 6).Value = 2.7648654815859
$ws.Cells.Item(13, 7).Value = 0.0024917316166293
$ws.Cells.Item(13, 9).Value = 1.222498342292681
$ws.Cells.Item(13, 10).Value = 0.2135327036101984
$ws.Cells.Item(13, 13).Value = 0.743403876483697
$ws.Cells.Item(13, 14).Value = 1.645017981342271

$ws.Cells.Item(14, 2).Value = 1.959830886587383
$ws.Cells.Item(14, 3).Value = 0.7481206758088774
$ws.Cells.Item(14, 4).Value = 0.09093736938120855
$ws.Cells.Item(14, 6).Value = 2.756659099121293
$ws.Cells.Item(14, 7).Value = 0.002492662355991656
$ws.Cells.Item(14, 9).Value = 1.220387884986565
$ws.Cells.Item(14, 10).Value = 0.2133795192768204
$ws.Cells.Item(14, 13).Value = 0.7370683121328483
$ws.Cells.Item(14, 14).Value = 1.648443208035957

$ws.Cells.Item(15, 2).Value = 1.945087599906969
$ws.Cells.Item(15, 3).Value = 0.74247262481299
$ws.Cells.Item(15, 4).Value = 0.09090304179759912
$ws.Cells.Item(15, 6).Value = 2.751665291174788
$ws.Cells.Item(15, 7).Value = 0.002493235703605266
$ws.Cells.Item(15, 9).Value = 1.219110037203194
$ws.Cells.Item(15, 10).Value = 0.2132875639656362
$ws.Cells.Item(15, 13).Value = 0.7331947779770616
$ws.Cells.Item(15, 14).Value = 1.650554334805406

$ws.Cells.Item(16, 2).Value = 1.860848026619863
$ws.Cells.Item(16, 3).Value = 0.7102135176241973
$ws.Cells.Item(16, 4).Value = 0.09071002589956478
$ws.Cells.Item(16, 6).Value = 2.723549727426388
$ws.Cells.Item(16, 7).Value = 0.00249657112815923
$ws.Cells.Item(16, 9).Value = 1.212014448391145
$ws.Cells.Item(16, 10).Value = 0.2127893074282952
$ws.Cells.Item(16, 13).Value = 0.7111069316819965
$ws.Cells.Item(16, 14).Value = 1.662852408372913

$ws.Cells.Item(17, 2).Value = 1.809390303511691
$ws.Cells.Item(17, 3).Value = 0.6905189109876346
$ws.Cells.Item(17, 4).Value = 0.09059493332641821
$ws.Cells.Item(17, 6).Value = 2.706748118083624
$ws.Cells.Item(17, 7).Value = 0.002498661856333187
$ws.Cells.Item(17, 9).Value = 1.207863602125869
$ws.Cells.Item(17, 10).Value = 0.2125092234010282
$ws.Cells.Item(17, 13).Value = 0.6976544526579431
$ws.Cells.Item(17, 14).Value = 1.670574997331983

$ws.Cells.Item(18, 2).Value = 1.779873838736989
$ws.Cells.Item(18, 3).Value = 0.6792258474801542
$ws.Cells.Item(18, 4).Value = 0.09052995926108665
$ws.Cells.Item(18, 6).Value = 2.697247726673936
$ws.Cells.Item(18, 7).Value = 0.002499880805929498
$ws.Cells.Item(18, 9).Value = 1.20555012704169
$ws.Cells.Item(18, 10).Value = 0.212357514867044
$ws.Cells.Item(18, 13).Value = 0.6899527337665177
$ws.Cells.Item(18, 14).Value = 1.675082127843531

$ws.Cells.Item(19, 2).Value = 1.769893887510705
$ws.Cells.Item(19, 3).Value = 0.6754081534137413
$ws.Cells.Item(19, 4).Value = 0.09050817091845076
$ws.Cells.Item(19, 6).Value = 2.694059074136419
$ws.Cells.Item(19, 7).Value = 0.002500296345615362
$ws.Cells.Item(19, 9).Value = 1.204779504406829
$ws.Cells.Item(19, 10).Value = 0.212307760277433
$ws.Cells.Item(19, 13).Value = 0.6873512006684592
$ws.Cells.Item(19, 14).Value = 1.676619365627381

$ws.Cells.Item(20, 2).Value = 1.81485970969311
$ws.Cells.Item(20, 3).Value = 0.6926118325224024
$ws.Cells.Item(20, 4).Value = 0.09060705855458195
$ws.Cells.Item(20, 6).Value = 2.708519748332264
$ws.Cells.Item(20, 7).Value = 0.002498437596566205
$ws.Cells.Item(20, 9).Value = 1.208297803723028
$ws.Cells.Item(20, 10).Value = 0.2125380668874044
$ws.Cells.Item(20, 13).Value = 0.6990827852327897
$ws.Cells.Item(20, 14).Value = 1.6697461525739

$ws.Cells.Item(21, 2).Value = 1.96690509838794
$ws.Cells.Item(21, 3).Value = 0.7508309827408652
$ws.Cells.Item(21, 4).Value = 0.09095389568091861
$ws.Cells.Item(21, 6).Value = 2.759062665768056
$ws.Cells.Item(21, 7).Value = 0.00249238829997586
$ws.Cells.Item(21, 9).Value = 1.221004677363041
$ws.Cells.Item(21, 10).Value = 0.2134241225112845
$ws.Cells.Item(21, 13).Value = 0.738927723603183
$ws.Cells.Item(21, 14).Value = 1.647434409906296

$ws.Cells.Item(22, 2).Value = 2.067002074598349
$ws.Cells.Item(22, 3).Value = 0.7891958762147624
$ws.Cells.Item(22, 4).Value = 0.09119136153979568
$ws.Cells.Item(22, 6).Value = 2.793563547291342
$ws.Cells.Item(22, 7).Value = 0.002488579982937965
$ws.Cells.Item(22, 9).Value = 1.229974073527487
$ws.Cells.Item(22, 10).Value = 0.2140870767658427
$ws.Cells.Item(22, 13).Value = 0.7652898352189652
$ws.Cells.Item(22, 14).Value = 1.633437247335472

$ws.Cells.Item(23, 2).Value = 2.013511549362818
$ws.Cells.Item(23, 3).Value = 0.7686906735322623
$ws.Cells.Item(23, 4).Value = 0.09106363986512989
$ws.Cells.Item(23, 6).Value = 2.775014821737898
$ws.Cells.Item(23, 7).Value = 0.0024905993093301
$ws.Cells.Item(23, 9).Value = 1.225125827147792
$ws.Cells.Item(23, 10).Value = 0.2137255589156766
$ws.Cells.Item(23, 13).Value = 0.751190376213799
$ws.Cells.Item(23, 14).Value = 1.640854118701235

$ws.Cells.Item(24, 2).Value = 1.812386781730538
$ws.Cells.Item(24, 3).Value = 0.6916655304755182
$ws.Cells.Item(24, 4).Value = 0.09060157301775718
$ws.Cells.Item(24, 6).Value = 2.707718298944997
$ws.Cells.Item(24, 7).Value = 0.002498538931554446
$ws.Cells.Item(24, 9).Value = 1.208101274173586
$ws.Cells.Item(24, 10).Value = 0.2125249977346897
$ws.Cells.Item(24, 13).Value = 0.6984369355244411
$ws.Cells.Item(24, 14).Value = 1.670120663666495

$ws.Cells.Item(25, 2).Value = 1.598519349572427
$ws.Cells.Item(25, 3).Value = 0.6099064769820757
$ws.Cells.Item(25, 4).Value = 0.09015002821531581
$ws.Cells.Item(25, 6).Value = 2.641369480317863
$ws.Cells.Item(25, 7).Value = 0.002507731389988577
$ws.Cells.Item(25, 9).Value = 1.192562485543945
$ws.Cells.Item(25, 10).Value = 0.2115887777677443
$ws.Cells.Item(25, 13).Value = 0.6429006001094706
$ws.Cells.Item(25, 14).Value = 1.704184144849677

Write-Host "Applied 216 cell updates"